# Make plot and visit forms more consistent:
#  - rename the "name" field value to "plot_name" on the survey and
#    settings sheets (adds a new shared string "plot_name")
#  - make "settings" the active/selected sheet instead of "survey"
#  - update the remembered cell selections on the affected sheets

$wb = $excel.ActiveWorkbook

$survey = $wb.Worksheets.Item("survey")
$settings = $wb.Worksheets.Item("settings")

# Update the "name" -> "plot_name" values
$survey.Range("D3").Value = "plot_name"
$settings.Range("B5").Value = "plot_name"

# Move the active selection on the survey sheet (it is no longer the
# active tab) and activate the settings sheet, which becomes selected.
[void]$survey.Range("D24").Select()

[void]$settings.Activate()
[void]$settings.Range("B5").Select()
